# ProductHierarchy.xlsx edit
# Commit: "New function; organized data after SQLite extraction; modified README.md"
#
# On the "Raw Data" sheet, the Tier value (column B) had been redundantly
# re-entered on several rows that belong to the same Product/Tier group as
# the row above, and row 10 repeated the Product label / Offer Number that
# is already implied by row 9's grouping. Clear the now-redundant cells so
# the sheet matches the cleaned-up hierarchy produced after the SQLite
# extraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw Data")

# Redundant "Tier" values within a Product/Tier group - clear them.
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B12").ClearContents()

# Row 10 repeated the Product label and Offer Number from its group header;
# clear them so only the Tier/Price/Probability remain for this sub-row.
$ws.Range("A10").ClearContents()
$ws.Range("C10").ClearContents()

# Reflect the final manual selection / active sheet left by the edit.
$ws.Activate()
$ws.Range("C10").Select() | Out-Null
